$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: RandomForest
$ws.Range("A2").Value = "RandomForest Multi-Max R"
$ws.Range("C2").Value = 0.9

# Row 3: XGBoost
$ws.Range("A3").Value = "XGBoost Multi-Max R"
$ws.Range("C3").Value = 0.95

# Row 4: Logistic Regression (label only, values unchanged)
$ws.Range("A4").Value = "Logistic Regression Multi-Max R"

# Row 5: Voting Classifier
$ws.Range("A5").Value = "Voting Classifier Multi-Max R"
$ws.Range("C5").Value = 0.95

# Row 6: Stacking Classifier
$ws.Range("A6").Value = "Stacking Classifier Multi-Max R"
$ws.Range("C6").Value = 0.9
